$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.81297300820363
$ws.Range("C2").Value = 6.798951297540322
$ws.Range("D2").Value = 5.992680736410231
$ws.Range("E2").Value = 16.39331199136112
$ws.Range("G2").Value = 37.48732519481784
$ws.Range("H2").Value = 16.11040189603196
$ws.Range("I2").Value = 23.09226271678976
$ws.Range("K2").Value = 11.09798392166926
$ws.Range("N2").Value = 19.04870061661193

$ws.Range("B3").Value = 10.49440811641334
$ws.Range("C3").Value = 6.48711000768777
$ws.Range("D3").Value = 5.875757886646279
$ws.Range("E3").Value = 15.46906176445426
$ws.Range("G3").Value = 37.21951205374302
$ws.Range("H3").Value = 16.12659173312948
$ws.Range("I3").Value = 23.11593736069569
$ws.Range("K3").Value = 10.86491339514205
$ws.Range("N3").Value = 19.11033559680503

$ws.Range("B4").Value = 10.29685371543023
$ws.Range("C4").Value = 6.289793328427747
$ws.Range("D4").Value = 5.804633534088312
$ws.Range("E4").Value = 14.87764814155789
$ws.Range("G4").Value = 37.06669540065129
$ws.Range("H4").Value = 16.13984090857926
$ws.Range("I4").Value = 23.13564811446016
$ws.Range("K4").Value = 10.72237897801617
$ws.Range("N4").Value = 19.15004717177155

$ws.Range("B5").Value = 10.21599552400927
$ws.Range("C5").Value = 6.208043256834748
$ws.Range("D5").Value = 5.775861394557248
$ws.Range("E5").Value = 14.63090054639666
$ws.Range("G5").Value = 37.00739471651742
$ws.Range("H5").Value = 16.14607001133255
$ws.Range("I5").Value = 23.14497744606305
$ws.Range("K5").Value = 10.66452761551736
$ws.Range("N5").Value = 19.16670037272604

$ws.Range("B6").Value = 10.20255194660634
$ws.Range("C6").Value = 6.194391746927129
$ws.Range("D6").Value = 5.771097872368463
$ws.Range("E6").Value = 14.5895905374883
$ws.Range("G6").Value = 36.99772883361545
$ws.Range("H6").Value = 16.14715441260957
$ws.Range("I6").Value = 23.1466047830868
$ws.Range("K6").Value = 10.65493810954236
$ws.Range("N6").Value = 19.16949406031146

$ws.Range("B7").Value = 10.29576447509335
$ws.Range("C7").Value = 6.28869607014607
$ws.Range("D7").Value = 5.804244587877641
$ws.Range("E7").Value = 14.87434326531738
$ws.Range("G7").Value = 37.06588354939084
$ws.Range("H7").Value = 16.13992155890139
$ws.Range("I7").Value = 23.13576868731515
$ws.Range("K7").Value = 10.72159771432086
$ws.Range("N7").Value = 19.15026985693504

$ws.Range("B8").Value = 10.70361787287913
$ws.Range("C8").Value = 6.692712955564215
$ws.Range("D8").Value = 5.952253617566519
$ws.Range("E8").Value = 16.07974341017384
$ws.Range("G8").Value = 37.39260358367928
$ws.Range("H8").Value = 16.11529633011725
$ws.Range("I8").Value = 23.09934944211724
$ws.Range("K8").Value = 11.01755527858507
$ws.Range("N8").Value = 19.06956526748762

$ws.Range("B9").Value = 11.48226368993143
$ws.Range("C9").Value = 7.4338333597952
$ws.Range("D9").Value = 6.245883139043928
$ws.Range("E9").Value = 18.29662607255325
$ws.Range("G9").Value = 38.1229852683037
$ws.Range("H9").Value = 16.09333844665327
$ws.Range("I9").Value = 23.06915407258315
$ws.Range("K9").Value = 11.59866124074989
$ws.Range("N9").Value = 18.92608257284422

$ws.Range("B10").Value = 12.03451956077778
$ws.Range("C10").Value = 7.941800452359926
$ws.Range("D10").Value = 6.461189625207735
$ws.Range("E10").Value = 19.92831104900199
$ws.Range("G10").Value = 38.71050563194082
$ws.Range("H10").Value = 16.09335449010456
$ws.Range("I10").Value = 23.07230006638408
$ws.Range("K10").Value = 12.0211481307089
$ws.Range("N10").Value = 18.82962015572745

$ws.Range("B11").Value = 12.28017938222956
$ws.Range("C11").Value = 8.164088586443228
$ws.Range("D11").Value = 6.558547521941488
$ws.Range("E11").Value = 20.62897916310933
$ws.Range("G11").Value = 38.98793727382948
$ws.Range("H11").Value = 16.09688291630052
$ws.Range("I11").Value = 23.07926377725567
$ws.Range("K11").Value = 12.21141520011849
$ws.Range("N11").Value = 18.78766976315046

$ws.Range("B12").Value = 12.37231109442952
$ws.Range("C12").Value = 8.246942079630017
$ws.Range("D12").Value = 6.595291043800813
$ws.Range("E12").Value = 20.88836008654002
$ws.Range("G12").Value = 39.09437376690623
$ws.Range("H12").Value = 16.09872599462523
$ws.Range("I12").Value = 23.08269796866794
$ws.Range("K12").Value = 12.28311286117767
$ws.Range("N12").Value = 18.77206112055939

$ws.Range("B13").Value = 12.35251005914966
$ws.Range("C13").Value = 8.22915773235953
$ws.Range("D13").Value = 6.587383781879199
$ws.Range("E13").Value = 20.83276171038927
$ws.Range("G13").Value = 39.07139077374705
$ws.Range("H13").Value = 16.0983065021911
$ws.Range("I13").Value = 23.08192287960147
$ws.Range("K13").Value = 12.26768831253062
$ws.Range("N13").Value = 18.77541040781484

$ws.Range("B14").Value = 12.287777544092
$ws.Range("C14").Value = 8.170931792918408
$ws.Range("D14").Value = 6.561573133644234
$ws.Range("E14").Value = 20.65043757693443
$ws.Range("G14").Value = 38.99666671218813
$ws.Range("H14").Value = 16.09702438530774
$ws.Range("I14").Value = 23.07953032632294
$ws.Range("K14").Value = 12.21732128811008
$ws.Range("N14").Value = 18.78638008132862

$ws.Range("B15").Value = 12.24800790341132
$ws.Range("C15").Value = 8.135092963469278
$ws.Range("D15").Value = 6.545746106692005
$ws.Range("E15").Value = 20.53798533281916
$ws.Range("G15").Value = 38.95107316753724
$ws.Range("H15").Value = 16.09630508273114
$ws.Range("I15").Value = 23.07816867323998
$ws.Range("K15").Value = 12.18642197975628
$ws.Range("N15").Value = 18.79313538741659

$ws.Range("B16").Value = 12.0183448825733
$ws.Range("C16").Value = 7.927091336914481
$ws.Range("D16").Value = 6.454811838202586
$ws.Range("E16").Value = 19.88168670095052
$ws.Range("G16").Value = 38.69257238708103
$ws.Range("H16").Value = 16.09319481236266
$ws.Range("I16").Value = 23.07195645166444
$ws.Range("K16").Value = 12.00866843710064
$ws.Range("N16").Value = 18.83240051483775

$ws.Range("B17").Value = 11.87595894604375
$ws.Range("C17").Value = 7.797193428517046
$ws.Range("D17").Value = 6.398847724157402
$ws.Range("E17").Value = 19.46844005232591
$ws.Range("G17").Value = 38.53653641791166
$ws.Range("H17").Value = 16.09218914015997
$ws.Range("I17").Value = 23.06956380717623
$ws.Range("K17").Value = 11.89907425826928
$ws.Range("N17").Value = 18.85698250463791

$ws.Range("B18").Value = 11.79354362743143
$ws.Range("C18").Value = 7.721655188270116
$ws.Range("D18").Value = 6.366605717170253
$ws.Range("E18").Value = 19.22683449313707
$ws.Range("G18").Value = 38.44775066817573
$ws.Range("H18").Value = 16.09194211339759
$ws.Range("I18").Value = 23.06870831448965
$ws.Range("K18").Value = 11.83586169603572
$ws.Range("N18").Value = 18.87130318914718

$ws.Range("B19").Value = 11.76555315407794
$ws.Range("C19").Value = 7.695939585690023
$ws.Range("D19").Value = 6.355681275008819
$ws.Range("E19").Value = 19.14435720511515
$ws.Range("G19").Value = 38.41785703164216
$ws.Range("H19").Value = 16.09191537048263
$ws.Range("I19").Value = 23.0685080304729
$ws.Range("K19").Value = 11.81443100761279
$ws.Range("N19").Value = 18.87618316517464

$ws.Range("B20").Value = 11.89117056695292
$ws.Range("C20").Value = 7.81110705828709
$ws.Range("D20").Value = 6.404810970423061
$ws.Range("E20").Value = 19.51283601051215
$ws.Range("G20").Value = 38.55304773775816
$ws.Range("H20").Value = 16.09226189082901
$ws.Range("I20").Value = 23.06976460525857
$ws.Range("K20").Value = 11.91075961567787
$ws.Range("N20").Value = 18.8543469014163

$ws.Range("B21").Value = 12.30681602970782
$ws.Range("C21").Value = 8.188070465377125
$ws.Range("D21").Value = 6.569158011144305
$ws.Range("E21").Value = 20.70415171959583
$ws.Range("G21").Value = 39.01857822442889
$ws.Range("H21").Value = 16.09738721267052
$ws.Range("I21").Value = 23.08021143188756
$ws.Range("K21").Value = 12.23212544641744
$ws.Range("N21").Value = 18.78315050755401

$ws.Range("B22").Value = 12.57321427599325
$ws.Range("C22").Value = 8.426710393017915
$ws.Range("D22").Value = 6.675829876719695
$ws.Range("E22").Value = 21.44809814087197
$ws.Range("G22").Value = 39.33082771459422
$ws.Range("H22").Value = 16.10369179408878
$ws.Range("I22").Value = 23.09168597934714
$ws.Range("K22").Value = 12.44007107194454
$ws.Range("N22").Value = 18.73823421327259

$ws.Range("B23").Value = 12.43154165233445
$ws.Range("C23").Value = 8.300067737916541
$ws.Range("D23").Value = 6.618977167625189
$ws.Range("E23").Value = 21.05419890833407
$ws.Range("G23").Value = 39.16347049169655
$ws.Range("H23").Value = 16.1000564310307
$ws.Range("I23").Value = 23.08513622059663
$ws.Range("K23").Value = 12.32930118001497
$ws.Range("N23").Value = 18.76205933208676

$ws.Range("B24").Value = 11.88429511689802
$ws.Range("C24").Value = 7.804819374863632
$ws.Range("D24").Value = 6.40211519400296
$ws.Range("E24").Value = 19.49277713640977
$ws.Range("G24").Value = 38.54558009365774
$ws.Range("H24").Value = 16.09222796877057
$ws.Range("I24").Value = 23.06967220455939
$ws.Range("K24").Value = 11.90547729957787
$ws.Range("N24").Value = 18.85553787141644

$ws.Range("B25").Value = 11.27463513119516
$ws.Range("C25").Value = 7.239422924062872
$ws.Range("D25").Value = 6.16634793271312
$ws.Range("E25").Value = 17.68113758635671
$ws.Range("G25").Value = 37.9161651726831
$ws.Range("H25").Value = 16.09644855757569
$ws.Range("I25").Value = 23.07288682364114
$ws.Range("K25").Value = 11.44189689975584
$ws.Range("N25").Value = 18.96332207094708
